# Apply crypto price/volume updates per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.766.77"
$ws.Range("E2").Value = "  +3.41%  "

$ws.Range("D3").Value = "2.219.41"
$ws.Range("E3").Value = "  +2.76%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.17"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("E6").Value = "  +1.69%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "64.88"
$ws.Range("E7").Value = "  +0.84%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.407"
$ws.Range("E9").Value = "  +2.29%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0881"
$ws.Range("E10").Value = "  +2.14%  "

$ws.Range("D12").Value = "2.549.34"
$ws.Range("E12").Value = "  +2.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "16.05"
$ws.Range("E13").Value = "  +0.37%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.37"
$ws.Range("E14").Value = "  +0.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.825"
$ws.Range("E15").Value = "  +1.32%  "

$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("D17").Value = "2.223.27"
$ws.Range("E17").Value = "  +3.03%  "

$ws.Range("D18").Value = "40.666.07"
$ws.Range("E18").Value = "  +3.31%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.88"
$ws.Range("E19").Value = "  +2.78%  "

$ws.Range("D20").Value = "0.0₃0904"
$ws.Range("E20").Value = "  +5.53%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.16"
$ws.Range("E21").Value = "  +0.71%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "253.31"
$ws.Range("E22").Value = "  +9.37%  "

$ws.Range("E24").Value = "  +1.24%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.30"
$ws.Range("E25").Value = "  -8.40%  "

$ws.Range("E26").Value = "  +2.27%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "173.06"
$ws.Range("E27").Value = "  +0.45%  "

$ws.Range("E28").Value = "  +4.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.32"
$ws.Range("E29").Value = "  +1.92%  "

$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.85"
$ws.Range("E31").Value = "  +6.87%  "

$ws.Range("E32").Value = "  +1.38%  "

$ws.Range("E33").Value = "  +0.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.21"
$ws.Range("E34").Value = "  +2.09%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.81"
$ws.Range("E35").Value = "  +1.14%  "

$ws.Range("E36").Value = "  +1.91%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +6.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.47"
$ws.Range("E38").Value = "  +2.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.89"
$ws.Range("E40").Value = "  +13.46%  "

$ws.Range("E41").Value = "  +1.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.73"
$ws.Range("E42").Value = "  +11.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.87"
$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("B44").Value = "TerraClassic"
$ws.Range("C44").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.000220"
$ws.Range("E44").Value = "  +45.60%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.24"
$ws.Range("E45").Value = "  +4.60%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.48"
$ws.Range("E46").Value = "  -1.90%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.512.77"
$ws.Range("E47").Value = "  -1.79%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0942"
$ws.Range("E48").Value = "  +1.85%  "

$ws.Range("B49").Value = "HuobiToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.86"
$ws.Range("E49").Value = "  +1.38%  "

$ws.Range("E50").Value = "  +1.77%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.73"
$ws.Range("E51").Value = "  +11.65%  "
